$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("RUNMANAGER")
$ws2 = $wb.Worksheets.Item("DATAMANAGER")

# --- RUNMANAGER: update Count for test1 row from 2 to 1 ---
$ws1.Range("D2").Value = "'1"

# --- DATAMANAGER: insert url/browser columns between Execute and uname ---
$ws2.Columns("C:D").Insert()

$url = "https://opensource-demo.orangehrmlive.com/web/index.php/auth/login"

$ws2.Range("D1").Value = "browser"
$ws2.Range("D2").Value = "chrome"
$ws2.Range("D3").Value = "firefox"
$ws2.Range("D4").Value = "chrome"
$ws2.Range("D5").Value = "firefox"

$ws2.Range("C1").Value = "url"
$ws2.Range("C2").Value = $url
$ws2.Range("C3").Value = $url
$ws2.Range("C4").Value = $url
$ws2.Range("C5").Value = $url

$ws2.Columns("C").ColumnWidth = 67.8

# Update selections on both sheets; re-select RUNMANAGER last so it stays
# the active tab (tabSelected) as in the target workbook.
$ws2.Range("B4").Select()
$ws1.Range("C3").Select()
